$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert a new blank column before N. This shifts N..T to O..U.
$ws.Columns("N:N").Insert()

# 2. Fill the new column N with header "nomePlano" and value "100MB_SPEEDY".
#    The inserted column already inherits the header style (bold) from its
#    left neighbor, so we only need to set the values.
$ws.Range("N1").Value = "nomePlano"
$ws.Range("N2").Value = "100MB_SPEEDY"

# 3. idContratoIXC shifted to O; update its value 154045 -> 154047
#    (leading apostrophe forces Excel to store the numeric-looking text as
#    text, matching the source data's string typing)
$ws.Range("O2").Value = "'154047"

# 4. Delete the old "logRetorno" column (now at Q)
$ws.Columns("Q:Q").Delete()

# 5. Delete the old "idGrupo" column (now at Q again after previous delete)
$ws.Columns("Q:Q").Delete()

# 6. Rename idRadusuarios header (now at Q) to idRadUsuarios
$ws.Range("Q1").Value = "idRadUsuarios"

# 7. Update idRadUsuarios value 181280 -> 181284
$ws.Range("Q2").Value = "'181284"

# 8. Update idClienteIXC value 117695 -> 117696
$ws.Range("C2").Value = "'117696"

# 9. Update logRetornoRad (now at S2) text with new id and id_cliente values
$ws.Range("S2").Value = "{'type': 'success', 'message': 'Registro inserido com sucesso!', 'id': '181284', 'atualiza_campos': [{'tipo': 's', 'campo': 'online', 'valor': ''}, {'tipo': 'i', 'campo': 'id_cliente', 'valor': '117696'}]}"
